$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Fecha) date shuffles
$ws.Range("D2").Value = 44330
$ws.Range("D3").Value = 44302
$ws.Range("D4").Value = 44313
$ws.Range("D5").Value = 44322
$ws.Range("D6").Value = 44327
$ws.Range("D7").Value = 44309
$ws.Range("D9").Value = 44306
$ws.Range("D10").Value = 44323

# Column M (Volumen) changes
$ws.Range("M4").Value = 120
$ws.Range("M6").Value = 60
$ws.Range("M7").Value = 80
$ws.Range("M10").Value = 80

# Row 6 Q/S/T changes (swap values between row6 and row7)
$ws.Range("Q6").Value = "$/caja 10 kilos empedrada"
$ws.Range("S6").Value = 11500
$ws.Range("T6").Value = 1

# Row 7 Q/S/T changes
$ws.Range("Q7").Value = "$/caja 14 kilos granel"
$ws.Range("S7").Value = 821
$ws.Range("T7").Value = 14
